{"js": "// Replace the multiplication expressions in the practice-sheet table with\n// the newly generated problems. Each old expression is unique within the\n// document, so a targeted search + replace per pair is safe and precise.\nconst replacements = [\n  [\"57\u00d744=\", \"20\u00d785=\"],\n  [\"38\u00d724=\", \"27\u00d793=\"],\n  [\"66\u00d736=\", \"28\u00d712=\"],\n  [\"33\u00d770=\", \"33\u00d727=\"],\n  [\"36\u00d735=\", \"64\u00d785=\"],\n  [\"40\u00d726=\", \"81\u00d745=\"],\n  [\"95\u00d785=\", \"99\u00d755=\"],\n  [\"83\u00d799=\", \"55\u00d719=\"],\n  [\"54\u00d768=\", \"70\u00d784=\"],\n  [\"91\u00d767=\", \"12\u00d782=\"],\n  [\"30\u00d764=\", \"87\u00d736=\"],\n  [\"96\u00d781=\", \"69\u00d741=\"],\n  [\"44\u00d711=\", \"44\u00d712=\"],\n  [\"89\u00d757=\", \"41\u00d765=\"],\n  [\"20\u00d766=\", \"52\u00d786=\"],\n  [\"98\u00d774=\", \"37\u00d760=\"],\n  [\"62\u00d751=\", \"45\u00d742=\"],\n  [\"13\u00d717=\", \"11\u00d737=\"],\n  [\"77\u00d739=\", \"30\u00d751=\"],\n  [\"70\u00d749=\", \"95\u00d772=\"],\n  [\"50\u00d787=\", \"14\u00d725=\"],\n  [\"57\u00d735=\", \"85\u00d769=\"],\n  [\"19\u00d717=\", \"88\u00d751=\"],\n  [\"37\u00d730=\", \"82\u00d759=\"],\n  [\"85\u00d737=\", \"75\u00d729=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication expressions in the practice-sheet table with\n# the newly generated problems. Each old expression is unique within the\n# document, so Find/Replace per pair (re-scoped to the full document range\n# each time) is safe and precise.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"57\u00d744=\", \"20\u00d785=\"),\n    @(\"38\u00d724=\", \"27\u00d793=\"),\n    @(\"66\u00d736=\", \"28\u00d712=\"),\n    @(\"33\u00d770=\", \"33\u00d727=\"),\n    @(\"36\u00d735=\", \"64\u00d785=\"),\n    @(\"40\u00d726=\", \"81\u00d745=\"),\n    @(\"95\u00d785=\", \"99\u00d755=\"),\n    @(\"83\u00d799=\", \"55\u00d719=\"),\n    @(\"54\u00d768=\", \"70\u00d784=\"),\n    @(\"91\u00d767=\", \"12\u00d782=\"),\n    @(\"30\u00d764=\", \"87\u00d736=\"),\n    @(\"96\u00d781=\", \"69\u00d741=\"),\n    @(\"44\u00d711=\", \"44\u00d712=\"),\n    @(\"89\u00d757=\", \"41\u00d765=\"),\n    @(\"20\u00d766=\", \"52\u00d786=\"),\n    @(\"98\u00d774=\", \"37\u00d760=\"),\n    @(\"62\u00d751=\", \"45\u00d742=\"),\n    @(\"13\u00d717=\", \"11\u00d737=\"),\n    @(\"77\u00d739=\", \"30\u00d751=\"),\n    @(\"70\u00d749=\", \"95\u00d772=\"),\n    @(\"50\u00d787=\", \"14\u00d725=\"),\n    @(\"57\u00d735=\", \"85\u00d769=\"),\n    @(\"19\u00d717=\", \"88\u00d751=\"),\n    @(\"37\u00d730=\", \"82\u00d759=\"),\n    @(\"85\u00d737=\", \"75\u00d729=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
